# Plots.xlsx update — Importer config grouping (#779)
# "Fix Plots.xlsx data names" / "Fixed observed data naming (MALE)"
#
# The observed-data name used on the DataCombined sheet is updated to
# include the subject's gender ("MALE") in the naming pattern, matching
# the renamed Importer configuration pattern.

$wb = $excel.ActiveWorkbook

$oldName = "Laskin 1982.Group A_Aciclovir_1_Human_PeripheralVenousBlood_Plasma_2.5 mg/kg_iv_"
$newName = "Laskin 1982.Group A_Aciclovir_1_Human_MALE_PeripheralVenousBlood_Plasma_2.5 mg/kg_iv_"

$dataCombined = $wb.Worksheets.Item("DataCombined")

# Update the two "dataSet" cells (column F) that reference the observed
# data name, on the individual (row 3) and population (row 5) entries.
for ($row = 2; $row -le $dataCombined.UsedRange.Rows.Count; $row++) {
    $cell = $dataCombined.Cells.Item($row, 6)
    if ($cell.Value() -eq $oldName) {
        $cell.Value = $newName
    }
}

# Sheet/selection state: DataCombined becomes the selected/active tab with
# F6 selected, and plotConfiguration (previously the active tab) is left
# selected at M2 but is no longer the active sheet.
$plotConfiguration = $wb.Worksheets.Item("plotConfiguration")
[void]$plotConfiguration.Activate()
[void]$plotConfiguration.Range("M2").Select()

[void]$dataCombined.Activate()
[void]$dataCombined.Range("F6").Select()
